$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the OAMP model name (D7): OPA336N -> MAX 4238/4239
$ws.Range("D7").Value = "MAX 4238/4239"

# Update Precio Un.(USD) for OAMP (E7): 1.29 -> 1.52
$ws.Range("E7").Value = 1.52

# Update Cantidad for OAMP (G7): 1000 -> 2000
$ws.Range("G7").Value = 2000

# Update the selected cell to K6, matching the saved sheet view
$ws.Range("K6").Select()

$wb.Save()
